# Insert a new weekly data row for "Jengibre" (Vega Modelo de Temuco) at row 128.
# Excel's native row Insert() shifts the existing rows 128-183 down to 129-184,
# carrying their formatting (e.g. the date-format style on column D) with them -
# matching the rest of the diff exactly. We then populate the freshly inserted
# row 128 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(128).Insert()

$ws.Cells.Item(128, 1).Value = 10
$ws.Cells.Item(128, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(128, 3).Value = 'La Araucanía'
$ws.Cells.Item(128, 4).Value = 44784
$ws.Cells.Item(128, 5).Value = 9
$ws.Cells.Item(128, 6).Value = 100114007
$ws.Cells.Item(128, 7).Value = 'Jengibre'
$ws.Cells.Item(128, 8).Value = 'Sin especificar'
$ws.Cells.Item(128, 9).Value = 'Primera'
$ws.Cells.Item(128, 10).Value = 100
$ws.Cells.Item(128, 11).Value = 16000
$ws.Cells.Item(128, 12).Value = 16000
$ws.Cells.Item(128, 13).Value = 16000
$ws.Cells.Item(128, 14).Value = '$/caja 13 kilos'
$ws.Cells.Item(128, 15).Value = 'Perú'
$ws.Cells.Item(128, 16).Value = 1231
$ws.Cells.Item(128, 17).Value = 13
$ws.Cells.Item(128, 18).Value = 'Hortaliza'
